$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errores")

# Make "Errores" the active sheet (moves tabSelected/activeTab here).
$ws.Activate()

# Insert two new rows before row 15, pushing the existing rows 15 and 17
# down to rows 17 and 19 respectively.
$ws.Rows.Item(14).Insert()
$ws.Rows.Item(14).Insert()

# Fill in the two new comment rows.
$ws.Range("C14").Value = "Los campos de fecha en la base de datos estan como Date, lo cual impide que se almacene la hora, se cambian a datetime"
$ws.Range("C15").Value = "El mapeo de hibernate tiene los campos fecha como DATE se cambian a TIMESTAMP"

# Update the selection to match the new active cell.
$ws.Range("C16").Select()
